$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 15775.556
$ws.Range("I51").Value = 12995
$ws.Range("J51").Value = 18000
$ws.Range("K51").Value = 12995
$ws.Range("L51").Value = 18000
$ws.Range("M51").Value = -12511

$ws.Range("H70").Value = 2957.7896
$ws.Range("I70").Value = 2626.7693
$ws.Range("J70").Value = 3675
$ws.Range("K70").Value = 7880.3079
$ws.Range("L70").Value = 11025
$ws.Range("M70").Value = -7610.3079
$ws.Range("N70").Value = -11565

$ws.Range("H73").Value = 2957.7896
$ws.Range("I73").Value = 2626.7693
$ws.Range("J73").Value = 3675
$ws.Range("K73").Value = 7880.3079
$ws.Range("L73").Value = 11025
$ws.Range("M73").Value = -6944.3079
$ws.Range("N73").Value = -12897

$ws.Range("H103").Value = 1004
$ws.Range("I103").Value = 1004
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 3012
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -2426
$ws.Range("N103").ClearContents()

$ws.Range("H116").Value = 4750
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 4750
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 4750
$ws.Range("N116").Value = -11634
$ws.Range("M116").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H38").Value = 2541.625
$ws.Range("I38").Value = 2044.5714
$ws.Range("J38").Value = 6021
$ws.Range("K38").Value = 2044.5714
$ws.Range("L38").Value = 6021
$ws.Range("M38").Value = -1577.5714

$ws.Range("H41").Value = 568.3333
$ws.Range("I41").Value = 568.3333
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 568.3333
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -154.3333

$ws.Range("H61").Value = 5598.5
$ws.Range("I61").Value = 5598.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 5598.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -5386.5
$ws.Range("N61").ClearContents()

$ws.Range("H63").Value = 4248.8887
$ws.Range("I63").Value = 1373.3334
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 1373.3334
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = -687.3334

$ws.Range("H66").Value = 4248.8887
$ws.Range("I66").Value = 1373.3334
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 6866.666999999999
$ws.Range("L66").Value = 50000
$ws.Range("M66").Value = -3434.666999999999

$ws.Range("H136").Value = 5598.5
$ws.Range("I136").Value = 5598.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 16795.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -14245.5
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 5856.75
$ws.Range("I75").Value = 5856.75
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 5856.75
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -4920.75

$ws.Range("H78").Value = 5856.75
$ws.Range("I78").Value = 5856.75
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 17570.25
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -12890.25

$ws.Range("H82").Value = 41277
$ws.Range("I82").Value = 20395.857
$ws.Range("J82").Value = 89999.664
$ws.Range("K82").Value = 20395.857
$ws.Range("L82").Value = 89999.664
$ws.Range("M82").Value = -20012.857

$ws.Range("H85").Value = 41277
$ws.Range("I85").Value = 20395.857
$ws.Range("J85").Value = 89999.664
$ws.Range("K85").Value = 20395.857
$ws.Range("L85").Value = 89999.664
$ws.Range("M85").Value = -19069.857

$ws.Range("H102").Value = 20430.285
$ws.Range("I102").Value = 20430.285
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 20430.285
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -17185.285
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 340.5
$ws.Range("I16").Value = 263.91666
$ws.Range("J16").Value = 800
$ws.Range("K16").Value = 263.91666
$ws.Range("L16").Value = 800
$ws.Range("M16").Value = 23.08334000000002
$ws.Range("N16").Value = -1374

$ws.Range("H33").Value = 16859.75
$ws.Range("I33").Value = 2479.8333
$ws.Range("J33").Value = 59999.5
$ws.Range("K33").Value = 2479.8333
$ws.Range("L33").Value = 59999.5
$ws.Range("M33").Value = -2100.8333

$ws.Range("H36").Value = 36666.168
$ws.Range("I36").Value = 19999.25
$ws.Range("J36").Value = 70000
$ws.Range("K36").Value = 19999.25
$ws.Range("L36").Value = 70000
$ws.Range("M36").Value = -19611.25

$ws.Range("H40").Value = 36666.168
$ws.Range("I40").Value = 19999.25
$ws.Range("J40").Value = 70000
$ws.Range("K40").Value = 19999.25
$ws.Range("L40").Value = 70000
$ws.Range("M40").Value = -19839.25

$ws.Range("H68").Value = 90000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 90000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 90000
$ws.Range("N68").Value = -91498

$ws.Range("H71").Value = 90000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 90000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 270000
$ws.Range("N71").Value = -277488

$ws.Range("H94").Value = 1777.25
$ws.Range("I94").Value = 1312
$ws.Range("J94").Value = 1932.3334
$ws.Range("K94").Value = 1312
$ws.Range("L94").Value = 1932.3334
$ws.Range("M94").Value = -861
$ws.Range("N94").Value = -2834.3334

$ws.Range("H113").Value = 340.5
$ws.Range("I113").Value = 263.91666
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 263.91666
$ws.Range("L113").Value = 800
$ws.Range("M113").Value = 1906.08334
$ws.Range("N113").Value = -5140

$ws.Range("H134").Value = 3324.3635
$ws.Range("I134").Value = 1801.4445
$ws.Range("J134").Value = 10177.5
$ws.Range("K134").Value = 5404.333500000001
$ws.Range("L134").Value = 30532.5
$ws.Range("M134").Value = -2869.333500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 899.25
$ws.Range("I6").Value = 48.5
$ws.Range("J6").Value = 1750
$ws.Range("K6").Value = 145.5
$ws.Range("L6").Value = 5250
$ws.Range("M6").Value = -32.5

$ws.Range("H32").Value = 8750
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 8750
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 26250
$ws.Range("N32").Value = -26816

$ws.Range("H34").Value = 4216.5
$ws.Range("I34").Value = 149.5
$ws.Range("J34").Value = 6250
$ws.Range("K34").Value = 448.5
$ws.Range("L34").Value = 18750
$ws.Range("M34").Value = -364.5
$ws.Range("N34").Value = -18918

$ws.Range("H74").Value = 2555
$ws.Range("I74").Value = 2555
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 7665
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -6604

$ws.Range("H75").Value = 2899.5
$ws.Range("I75").Value = 1000
$ws.Range("J75").Value = 3532.6667
$ws.Range("K75").Value = 3000
$ws.Range("L75").Value = 10598.0001
$ws.Range("M75").Value = -2002
$ws.Range("N75").Value = -12594.0001

$ws.Range("H77").Value = 2555
$ws.Range("I77").Value = 2555
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 22995
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -17691

$ws.Range("H78").Value = 2899.5
$ws.Range("I78").Value = 1000
$ws.Range("J78").Value = 3532.6667
$ws.Range("K78").Value = 9000
$ws.Range("L78").Value = 31794.0003
$ws.Range("M78").Value = -4008
$ws.Range("N78").Value = -41778.0003

$ws.Range("H80").Value = 999
$ws.Range("I80").Value = 999
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2997
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2061
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 999
$ws.Range("I83").Value = 999
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 8991
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -4311
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 4001.5
$ws.Range("I12").Value = 4001.5
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 4001.5
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -3861.5

$ws.Range("H26").Value = 9642.857
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 9642.857
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 9642.857
$ws.Range("N26").Value = -10202.857

$ws.Range("H38").Value = 3000
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 3000
$ws.Range("N38").Value = -3926

$ws.Range("H50").Value = 9642.857
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 9642.857
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 9642.857
$ws.Range("N50").Value = -10638.857

$ws.Range("H124").Value = 100000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 100000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 100000
$ws.Range("N124").Value = -109820

$ws.Range("H126").Value = 1949.25
$ws.Range("I126").Value = 1965.6666
$ws.Range("J126").Value = 1900
$ws.Range("K126").Value = 5896.9998
$ws.Range("L126").Value = 5700
$ws.Range("M126").Value = -3426.9998
$ws.Range("N126").Value = -10640

$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -205

$ws.Range("H27").Value = 500
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -393

$ws.Range("H46").Value = 3998.75
$ws.Range("I46").Value = 1995
$ws.Range("J46").Value = 4666.6665
$ws.Range("K46").Value = 1995
$ws.Range("L46").Value = 4666.6665
$ws.Range("M46").Value = -1807
$ws.Range("N46").Value = -5042.6665

$ws.Range("H132").Value = 3369.8
$ws.Range("I132").Value = 3369.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10109.4
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7579.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H29").Value = 4999.5
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 4999.5
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 4999.5
$ws.Range("N29").Value = -5579.5
$ws.Range("M29").ClearContents()

$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H104").Value = 23331.666
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 23331.666
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 23331.666
$ws.Range("N104").Value = -30319.666
